$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header: "ID-Number" -> "ID-Number / Topic"
$ws.Range("D2").Value = "ID-Number / Topic"

# Insert a new row above the "Business rule Task with DMN" row (old row 17),
# pushing the rest of the process steps down by one row.
$ws.Rows.Item(17).Insert()

# The inserted row only carries a new variable name in column E.
$ws.Range("E17").Value = "email"

# The DMN gateway conditions now reference the renamed process variable
# "Output_absenceFromType" instead of "absence_from_type_decision", and the
# auto-approved branch condition literal changed from "approved" to
# "auto-approved".
$ws.Range("E22").Value = '${Output_absenceFromType=="rejected"}'
$ws.Range("E23").Value = '${Output_absenceFromType=="prov-approved"}'
$ws.Range("E24").Value = '${Output_absenceFromType=="auto-approved"}'

# New process step appended: a Service Task that informs the employee.
$ws.Range("B25").Value = "Service Task"
$ws.Range("C25").Value = "Inform the employee"
$ws.Range("D25").Value = "informEmployee"

# Grow Table1 so the new row is included.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:E35"))

# Restore the selection to where the user left off.
$ws.Range("E33").Select()
